$d = $word.ActiveDocument

# 1. Append a full-stop to the June 2nd diary entry (paragraph 4).
$p4 = $d.Paragraphs(4)
$p4.Range.InsertAfter("。")

# 2. Insert the June 3rd date line + entry, and the June 7th date line +
#    entry right after the (now-period-terminated) June 2nd entry.
$p4.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "2022年6月3日星期五"

$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "中雨，今天是农历五月初五，是中国的传统节日：端午节，这一天我们要吃粽子，赛龙舟。"

$d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = "2022年6月7日星期二"

$d.Paragraphs(7).Range.InsertParagraphAfter()
$d.Paragraphs(8).Range.Text = "晴，今天是高考的一天，上午考语文，下午考数学。今天天气不错，心情也很好"

# 3. After the bookmark paragraph (now paragraph 9), append two empty
#    paragraphs. InsertParagraphAfter() always seeds the new paragraph
#    with a (property-only) run, so create the two placeholders first
#    and then rewrite each placeholder's full range (mark included) with
#    bare OOXML so no stray <w:r> survives, matching the diff exactly.
$d.Content.InsertParagraphAfter()
$d.Content.InsertParagraphAfter()

$p10 = $d.Paragraphs(10)
$r10 = $d.Range($p10.Range.Start, $p10.Range.End)
$r10.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr></w:p>")

$p11 = $d.Paragraphs(11)
$r11 = $d.Range($p11.Range.Start, $p11.Range.End)
$r11.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:hint='default'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr></w:p>")

Write-Output ("final paragraph count: " + $d.Paragraphs.Count)
